# Weekly update: prepend a new price report row for
# "Comercializadora del Agro de Limarí - Haba" and push the existing
# history down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 47; this shifts rows
# 47:65 down to 48:66 and extends the used range to A1:R66.
$ws.Rows("47").Insert()

# Populate the newly inserted row 47 with this week's figures.
$ws.Cells.Item(47, 1).Value = 2
$ws.Cells.Item(47, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44798
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 100112026
$ws.Cells.Item(47, 7).Value = "Haba"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 1500
$ws.Cells.Item(47, 11).Value = 7500
$ws.Cells.Item(47, 12).Value = 8000
$ws.Cells.Item(47, 13).Value = 7750
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(47, 16).Value = 310
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
